$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1, matching the style of the existing header row (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# "Save" column values for rows 2-11
$saveValues = @(0, 0, 1, 0, 1, 0, 0, 0, 0, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
